$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: insert a new data row (ID=5, papaya/Fiji/pink) before the
# existing row 6 (banana), shifting the remaining rows down and carrying
# their formats with them. ---
$ws1.Rows.Item(6).Insert()
$ws1.Range("A6").Value = 44867
$ws1.Range("B6").Value = 5
$ws1.Range("C6").Value = "papaya"
$ws1.Range("D6").Value = "Fiji"
$ws1.Range("E6").Value = "pink"

# --- Sheet2: add the new row (ID=5, papaya/elongated/salmon) coming from
# the merged table, keeping the default (unstyled) cell format rather than
# copying the neighbouring row's style. ---
$ws2.Range("A8").Value = 5
$ws2.Range("B8").Value = "papaya"
$ws2.Range("C8").Value = "elongated"
$ws2.Range("D8").Value = "salmon"

# Restore each sheet's own selection as left in the source file. Select
# Sheet2's cell first so that selecting Sheet1's cell afterwards leaves
# Sheet1 as the active (tabSelected) sheet, matching the source workbook.
$ws2.Range("B9").Select()
$ws1.Range("A7").Select()
